$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.375.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.472.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.71%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.475.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.426"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.083.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("E15").Value = "  -2.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.378.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000175"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.480.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.63%  "

$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.534"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000121"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.56%  "

$ws.Range("E29").Value = "  -1.84%  "

$ws.Range("E30").Value = "  +0.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.43%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.95%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "24.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.76%  "

$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("E37").Value = "  -4.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.883"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.712.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0699"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0296"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "324.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.70%  "
